# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures on
# the cryptos worksheet with the latest scraped snapshot.
#
# Columns D/E are stored as plain TEXT, not numbers, even though most
# values look numeric (e.g. "0.999", "19.41"): many prices use "." as a
# thousands separator (e.g. "60.920.87") and several need to keep
# significant trailing zeros ("1.60", "0.0560", "61.10") that a real
# number would drop. Handing a numeric-looking string straight to
# Range.Value lets the COM layer auto-coerce it into a true number,
# silently destroying that formatting - so for anything that parses as a
# plain number we prefix it with a leading apostrophe first, exactly like
# typing '0.999 into a cell in the Excel UI to force text entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Looks-Numeric($Text) {
    return ($Text -match '^[+-]?(\d+\.?\d*|\.\d+)([eE][+-]?\d+)?$')
}

function Set-TextValue($CellRef, $Text) {
    if (Looks-Numeric $Text) {
        $ws.Range($CellRef).Value = "'" + $Text
    } else {
        $ws.Range($CellRef).Value = $Text
    }
}

Set-TextValue "D2" "60.920.87"
Set-TextValue "E2" "  +0.95%  "
Set-TextValue "D3" "2.601.99"
Set-TextValue "E3" "  +0.62%  "
Set-TextValue "E4" "  -0.16%  "
Set-TextValue "D5" "523.57"
Set-TextValue "E5" "  +3.25%  "
Set-TextValue "D6" "154.97"
Set-TextValue "E6" "  +1.15%  "
Set-TextValue "D7" "0.999"
Set-TextValue "E7" "  -0.07%  "
Set-TextValue "D8" "0.589"
Set-TextValue "E8" "  +1.66%  "
Set-TextValue "D9" "6.68"
Set-TextValue "E9" "  +2.28%  "
Set-TextValue "D10" "0.105"
Set-TextValue "E10" "  +1.39%  "
Set-TextValue "D11" "0.347"
Set-TextValue "E11" "  +0.06%  "
Set-TextValue "E12" "  +1.10%  "
Set-TextValue "D13" "3.058.14"
Set-TextValue "E13" "  +0.55%  "
Set-TextValue "D14" "60.942.54"
Set-TextValue "E14" "  +0.90%  "
Set-TextValue "D15" "21.68"
Set-TextValue "E15" "  +0.24%  "
Set-TextValue "D16" "0.0000141"
Set-TextValue "E16" "  +1.06%  "
Set-TextValue "D17" "2.607.97"
Set-TextValue "E17" "  +0.31%  "
Set-TextValue "D19" "355.43"
Set-TextValue "E19" "  +2.61%  "
Set-TextValue "D20" "10.57"
Set-TextValue "E20" "  +1.25%  "
Set-TextValue "D21" "6.21"
Set-TextValue "E21" "  +1.64%  "
Set-TextValue "D22" "0.998"
Set-TextValue "E22" "  +0.13%  "
Set-TextValue "D23" "61.10"
Set-TextValue "E23" "  +2.37%  "
Set-TextValue "D24" "0.426"
Set-TextValue "E24" "  +1.47%  "
Set-TextValue "E25" "  +0.13%  "
Set-TextValue "D26" "2.720.87"
Set-TextValue "E26" "  +1.06%  "
Set-TextValue "D27" "0.999"
Set-TextValue "E27" "  -0.09%  "
Set-TextValue "D28" "0.0₃0847"
Set-TextValue "E28" "  +0.16%  "
Set-TextValue "D29" "7.41"
Set-TextValue "E29" "  +0.45%  "
Set-TextValue "E30" "  -0.05%  "
Set-TextValue "D31" "6.26"
Set-TextValue "E31" "  +9.47%  "
Set-TextValue "D32" "19.41"
Set-TextValue "E32" "  +0.10%  "
Set-TextValue "D33" "1.60"
Set-TextValue "E33" "  +2.89%  "
Set-TextValue "D34" "148.45"
Set-TextValue "E34" "  -3.24%  "
Set-TextValue "D35" "4.19"
Set-TextValue "E35" "  +4.95%  "
Set-TextValue "D36" "1.20"
Set-TextValue "E36" "  +1.22%  "
Set-TextValue "D37" "0.909"
Set-TextValue "E37" "  +7.05%  "
Set-TextValue "D38" "0.894"
Set-TextValue "E38" "  +5.30%  "
Set-TextValue "D39" "1.50"
Set-TextValue "E39" "  +1.43%  "
Set-TextValue "E40" "  +1.58%  "
Set-TextValue "E41" "  +0.86%  "
Set-TextValue "D42" "291.26"
Set-TextValue "E42" "  -1.65%  "
Set-TextValue "D43" "0.101"
Set-TextValue "E43" "  +2.08%  "
Set-TextValue "D44" "0.622"
Set-TextValue "E44" "  +0.00%  "
Set-TextValue "D45" "0.0560"
Set-TextValue "E45" "  +0.40%  "
Set-TextValue "D46" "0.998"
Set-TextValue "E46" "  -0.05%  "
Set-TextValue "D47" "19.55"
Set-TextValue "E47" "  -1.33%  "
Set-TextValue "D48" "4.92"
Set-TextValue "E48" "  +0.73%  "
Set-TextValue "E49" "  +2.29%  "
Set-TextValue "E50" "  +0.29%  "
Set-TextValue "D51" "19.23"
Set-TextValue "E51" "  +8.97%  "
